$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the V+ input value from 4 to 3.3; dependent formulas (D5, D6, D7)
# will recalculate automatically.
$ws.Range("D3").Value = 3.3

# Update the active cell selection to C7 (was E7).
$ws.Range("C7").Select()

$excel.Calculate()
